# Auto-generated edit script applying scraped diff changes to Gungnir_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 26 (ALC)
$ws.Range("H26").Value = 10375
$ws.Range("I26").Value = 5500
$ws.Range("J26").Value = 12000
$ws.Range("K26").Value = 5500
$ws.Range("L26").Value = 12000
$ws.Range("M26").Value = -5156
$ws.Range("N26").Value = -12688

# Row 41 (ALC)
$ws.Range("H41").Value = 146.14285
$ws.Range("I41").Value = 95.53846
$ws.Range("J41").Value = 228.375
$ws.Range("K41").Value = 95.53846
$ws.Range("L41").Value = 228.375
$ws.Range("M41").Value = 344.46154

# Row 51 (ALC)
$ws.Range("H51").Value = 13290.412
$ws.Range("I51").Value = 50463
$ws.Range("J51").Value = 1852.6923
$ws.Range("K51").Value = 50463
$ws.Range("L51").Value = 1852.6923
$ws.Range("M51").Value = -49979
$ws.Range("N51").Value = -2820.6923

# Row 80 (ALC)
$ws.Range("H80").Value = 319.9565
$ws.Range("I80").Value = 286.5
$ws.Range("J80").Value = 440.4
$ws.Range("K80").Value = 859.5
$ws.Range("L80").Value = 1321.2
$ws.Range("M80").Value = 138.5
$ws.Range("N80").Value = -3317.2

# Row 83 (ALC)
$ws.Range("H83").Value = 319.9565
$ws.Range("I83").Value = 286.5
$ws.Range("J83").Value = 440.4
$ws.Range("K83").Value = 2578.5
$ws.Range("L83").Value = 3963.6
$ws.Range("M83").Value = 2413.5
$ws.Range("N83").Value = -13947.6

# Row 137 (ALC)
$ws.Range("H137").Value = 1068.5536
$ws.Range("I137").Value = 953.6739
$ws.Range("J137").Value = 1597
$ws.Range("K137").Value = 2861.0217
$ws.Range("L137").Value = 4791
$ws.Range("M137").Value = -311.0217000000002
$ws.Range("N137").Value = -9891

# Row 138 (ALC)
$ws.Range("H138").Value = 2004.7424
$ws.Range("I138").Value = 1038.9375
$ws.Range("J138").Value = 2913.7354
$ws.Range("K138").Value = 3116.8125
$ws.Range("L138").Value = 8741.206200000001
$ws.Range("M138").Value = 2023.1875
$ws.Range("N138").Value = -19021.2062

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 732.5700000000001
$ws.Range("I32").Value = 718.12634
$ws.Range("J32").Value = 1007
$ws.Range("K32").Value = 718.12634
$ws.Range("L32").Value = 1007
$ws.Range("M32").Value = -431.12634
$ws.Range("N32").Value = -1581

# Row 33 (ARM)
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()

# Row 61 (ARM)
$ws.Range("H61").Value = 1298.2632
$ws.Range("I61").Value = 952.70966
$ws.Range("J61").Value = 2828.5715
$ws.Range("K61").Value = 952.70966
$ws.Range("L61").Value = 2828.5715
$ws.Range("M61").Value = -740.70966

# Row 86 (ARM)
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()

# Row 89 (ARM)
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()

# Row 132 (ARM)
$ws.Range("H132").Value = 1963123.9
$ws.Range("I132").Value = 2063.182
$ws.Range("J132").Value = 7356041
$ws.Range("K132").Value = 6189.545999999999
$ws.Range("L132").Value = 22068123
$ws.Range("M132").Value = -3659.545999999999
$ws.Range("N132").Value = -22073183

# Row 136 (ARM)
$ws.Range("H136").Value = 1298.2632
$ws.Range("I136").Value = 952.70966
$ws.Range("J136").Value = 2828.5715
$ws.Range("K136").Value = 2858.12898
$ws.Range("L136").Value = 8485.7145
$ws.Range("M136").Value = -308.12898

$ws = $wb.Worksheets.Item("BSM")
# Row 37 (BSM)
$ws.Range("H37").Value = 1163
$ws.Range("I37").Value = 1576
$ws.Range("J37").Value = 750
$ws.Range("K37").Value = 1576
$ws.Range("L37").Value = 750
$ws.Range("M37").Value = -1439
$ws.Range("N37").Value = -1024

# Row 88 (BSM)
$ws.Range("H88").Value = 37571.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 37571.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 37571.5
$ws.Range("N88").Value = -38383.5

# Row 91 (BSM)
$ws.Range("H91").Value = 37571.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 37571.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 37571.5
$ws.Range("N91").Value = -40379.5

# Row 132 (BSM)
$ws.Range("H132").Value = 45608
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 45608
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 45608
$ws.Range("N132").Value = -55728

# Row 134 (BSM)
$ws.Range("H134").Value = 5297166
$ws.Range("I134").Value = 1380.7059
$ws.Range("J134").Value = 27804252
$ws.Range("K134").Value = 4142.1177
$ws.Range("L134").Value = 83412756
$ws.Range("M134").Value = -1607.1177
$ws.Range("N134").Value = -83417826

$ws = $wb.Worksheets.Item("CUL")
# Row 131 (CUL)
$ws.Range("H131").Value = 809.9299999999999
$ws.Range("I131").Value = 404.08334
$ws.Range("J131").Value = 865.2727
$ws.Range("K131").Value = 1212.25002
$ws.Range("L131").Value = 2595.8181
$ws.Range("M131").Value = 3827.74998
$ws.Range("N131").Value = -12675.8181

$ws = $wb.Worksheets.Item("GSM")
# Row 51 (GSM)
$ws.Range("H51").Value = 59866.668
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 59866.668
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 59866.668
$ws.Range("N51").Value = -60884.668

# Row 132 (GSM)
$ws.Range("H132").Value = 16360.467
$ws.Range("I132").Value = 21100
$ws.Range("J132").Value = 13990.7
$ws.Range("K132").Value = 63300
$ws.Range("L132").Value = 41972.10000000001
$ws.Range("M132").Value = -60770
$ws.Range("N132").Value = -47032.10000000001

$ws = $wb.Worksheets.Item("LTW")
# Row 46 (LTW)
$ws.Range("H46").Value = 1535.0454
$ws.Range("I46").Value = 657.2222
$ws.Range("J46").Value = 2142.7693
$ws.Range("K46").Value = 657.2222
$ws.Range("L46").Value = 2142.7693
$ws.Range("M46").Value = -469.2222
$ws.Range("N46").Value = -2518.7693

# Row 74 (LTW)
$ws.Range("H74").Value = 16315.4
$ws.Range("I74").Value = 10000
$ws.Range("J74").Value = 17894.25
$ws.Range("K74").Value = 10000
$ws.Range("L74").Value = 17894.25
$ws.Range("M74").Value = -9002
$ws.Range("N74").Value = -19890.25

# Row 77 (LTW)
$ws.Range("H77").Value = 16315.4
$ws.Range("I77").Value = 10000
$ws.Range("J77").Value = 17894.25
$ws.Range("K77").Value = 30000
$ws.Range("L77").Value = 53682.75
$ws.Range("M77").Value = -25008
$ws.Range("N77").Value = -63666.75

# Row 87 (LTW)
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

# Row 90 (LTW)
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

# Row 132 (LTW)
$ws.Range("H132").Value = 19236828
$ws.Range("I132").Value = 37039096
$ws.Range("J132").Value = 10378.32
$ws.Range("K132").Value = 111117288
$ws.Range("L132").Value = 31134.96
$ws.Range("M132").Value = -111114758
$ws.Range("N132").Value = -36194.96

# Row 136 (LTW)
$ws.Range("H136").Value = 36006340
$ws.Range("I136").Value = 12871692
$ws.Range("J136").Value = 250001860
$ws.Range("K136").Value = 38615076
$ws.Range("L136").Value = 750005580
$ws.Range("M136").Value = -38612526
$ws.Range("N136").Value = -750010680

$ws = $wb.Worksheets.Item("WVR")
# Row 122 (WVR)
$ws.Range("H122").Value = 18163.291
$ws.Range("I122").Value = 32748.625
$ws.Range("J122").Value = 2605.6
$ws.Range("K122").Value = 98245.875
$ws.Range("L122").Value = 7816.799999999999
$ws.Range("M122").Value = -95795.875
$ws.Range("N122").Value = -12716.8

# Row 132 (WVR)
$ws.Range("H132").Value = 19964.098
$ws.Range("I132").Value = 22187
$ws.Range("J132").Value = 10887.25
$ws.Range("K132").Value = 66561
$ws.Range("L132").Value = 32661.75
$ws.Range("M132").Value = -64031
$ws.Range("N132").Value = -37721.75
